# Update cryptos list data (prices and 1h volume %) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.658.12'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.599.67'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''211.37'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '''0.517'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -1.23%  '
$ws.Range("D10").Value = '''19.47'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D12").Value = '1.823.40'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '1.651.28'
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").Value = '''64.83'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '26.638.77'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '''208.26'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '''6.96'
$ws.Range("E21").Value = '  +2.96%  '
$ws.Range("D22").Value = '''4.27'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("D24").Value = '''8.85'
$ws.Range("E24").Value = '  -0.64%  '
$ws.Range("D25").Value = '''145.59'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").Value = '''15.31'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = '''1.16'
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '''0.658'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").Value = '1.281.13'
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").Value = '''2.44'
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("D37").Value = '''1.50'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("D39").Value = '''0.843'
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").Value = '''5.47'
$ws.Range("E41").Value = '  +1.03%  '
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").Value = '''63.95'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("D45").Value = '''0.920'
$ws.Range("E45").Value = '  +9.63%  '
$ws.Range("D46").Value = '1.736.45'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = '''89.82'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").Value = '''1.60'
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.102'
$ws.Range("E49").Value = '  +3.62%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0506'
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''7.45'
$ws.Range("E51").Value = '  -1.14%  '
